$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + date range) ---------------------
# These cells hold rich-text (multi-run) shared strings in the source file;
# only a portion of the text changes (the last run in each case). The COM
# object model here exposes the concatenated text, so we write the full
# updated string back.
$ws.Range("A8").Value = "Volume 32   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/26/2025  Through  6/1/2025"

# --- Weekly crime-stat table updates (rows 15-28) --------------------------
# Some cells flip between a literal "0"/"***.*" placeholder (shared string,
# style 13 = General "text" style) and a real numeric value (style 14 =
# "#,##0" integer, style 15 = "#,##0.0;-#,##0.0" percent-style). We set
# NumberFormat explicitly so the saved style index matches what real Excel
# would assign, then set the value. For numeric -> placeholder-text
# conversions we flip the cell to Text format, write the literal text, and
# copy the number format (General, style 13) from a cell that is already in
# that exact placeholder style (C23) so the destination style matches
# exactly rather than creating a new one.

# Row 15
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 2
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = 100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 23.076923076923
$ws.Range("L15").Value = 128.571428571429
$ws.Range("M15").Value = 77.777777777777
$ws.Range("N15").Value = 45.454545454545

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = -22.58064516129
$ws.Range("I16").Value = 120
$ws.Range("J16").Value = 155
$ws.Range("K16").Value = -22.58064516129
$ws.Range("L16").Value = 5.263157894736
$ws.Range("M16").Value = -14.893617021276
$ws.Range("N16").Value = -77.099236641221

# Row 17
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = -10.869565217391
$ws.Range("I17").Value = 223
$ws.Range("J17").Value = 230
$ws.Range("K17").Value = -3.043478260869
$ws.Range("L17").Value = 13.197969543147
$ws.Range("M17").Value = 75.590551181102
$ws.Range("N17").Value = 36.80981595092

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 67
$ws.Range("J18").Value = 94
$ws.Range("K18").Value = -28.723404255319
$ws.Range("L18").Value = 6.349206349206
$ws.Range("M18").Value = -51.449275362318
$ws.Range("N18").Value = -92.463442069741

# Row 19
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -8.928571428571
$ws.Range("I19").Value = 285
$ws.Range("J19").Value = 367
$ws.Range("K19").Value = -22.343324250681
$ws.Range("L19").Value = -12.307692307692
$ws.Range("M19").Value = 54.054054054054
$ws.Range("N19").Value = -47.416974169741

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 104
$ws.Range("J20").Value = 123
$ws.Range("K20").Value = -15.447154471544
$ws.Range("L20").Value = -18.75
$ws.Range("M20").Value = 6.122448979591
$ws.Range("N20").Value = -88.571428571428

# Row 21
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 158
$ws.Range("G21").Value = 180
$ws.Range("H21").Value = -12.222222222222
$ws.Range("I21").Value = 819
$ws.Range("J21").Value = 982
$ws.Range("K21").Value = -16.598778004073
$ws.Range("L21").Value = -1.798561151079
$ws.Range("M21").Value = 17.167381974248
$ws.Range("N21").Value = -73.112278397898

# Row 22
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -44.444444444444
$ws.Range("I22").Value = 26
$ws.Range("K22").Value = -16.129032258064
$ws.Range("L22").Value = -45.833333333333
$ws.Range("M22").Value = 73.333333333333

# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -13.888888888888
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 140
$ws.Range("H24").Value = -30
$ws.Range("I24").Value = 546
$ws.Range("J24").Value = 962
$ws.Range("K24").Value = -43.243243243243
$ws.Range("L24").Value = -36.214953271028
$ws.Range("M24").Value = 24.657534246575

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -47.058823529411
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 67
$ws.Range("H25").Value = -55.223880597014
$ws.Range("I25").Value = 223
$ws.Range("J25").Value = 539
$ws.Range("K25").Value = -58.627087198515
$ws.Range("L25").Value = -51.731601731601

# Row 26
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 30
$ws.Range("E26").Value = -23.333333333333
$ws.Range("F26").Value = 95
$ws.Range("G26").Value = 118
$ws.Range("H26").Value = -19.491525423728
$ws.Range("I26").Value = 459
$ws.Range("J26").Value = 524
$ws.Range("K26").Value = -12.404580152671
$ws.Range("L26").Value = 24.728260869565
$ws.Range("M26").Value = 19.843342036553

# Row 27
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 100
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = -12
$ws.Range("L27").Value = 15.78947368421

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 16
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 77.777777777777
$ws.Range("I28").Value = 53
$ws.Range("J28").Value = 59
$ws.Range("K28").Value = -10.169491525423
$ws.Range("L28").Value = -14.516129032258

